$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D5").Value = "western gulf stream index_2025-05-29.png"
$ws.Range("D5").Select()
